$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update step/result text to refer to "assessments" instead of "kpis"
$ws.Range("C2").Value = "Step 1: While logged out try to view assessments from the url"
$ws.Range("C3").Value = "Step 2: Log in and from main page click on Team Assessments"
$ws.Range("D3").Value = "I am shown a list of assessments of people that is on my team"
$ws.Range("C4").Value = "Step 3: While Logged in try to view another team's Assessments"

# Add Actual Result / Pass-Fail columns for the test rows
$ws.Range("F2").Value = "pass"
$ws.Range("F3").Value = "pass"
$ws.Range("E4").Value = "directed to another teams assessment"
$ws.Range("F4").Value = "fail"

# Scroll the view so column B is the first visible column, then
# leave F5 as the active/selected cell (matches the saved view state)
$win = $excel.ActiveWindow
$win.ScrollRow = 1
$win.ScrollColumn = 2
$ws.Range("F5").Select()
